$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 2.27
$ws.Range("I4").Value = 3.2
$ws.Range("T4").Value = 6.2
$ws.Range("U4").Value = 10
$ws.Range("V4").Value = 9.5
$ws.Range("W4").Value = 23
$ws.Range("X4").Value = 22
$ws.Range("AE4").Value = 7.5
$ws.Range("AF4").Value = 15
$ws.Range("AG4").Value = 12
$ws.Range("AH4").Value = 45
$ws.Range("AI4").Value = 35

# Row 6 updates
$ws.Range("G6").Value = 2.27
$ws.Range("H6").Value = 2.87
$ws.Range("I6").Value = 3.3
$ws.Range("M6").Value = 2.22
$ws.Range("N6").Value = 2.47
$ws.Range("O6").Value = 1.42
$ws.Range("R6").Value = 2.07
$ws.Range("T6").Value = 5.8
$ws.Range("U6").Value = 9.5
$ws.Range("V6").Value = 9.75
$ws.Range("W6").Value = 23
$ws.Range("Z6").Value = 6.1
$ws.Range("AA6").Value = 5.8
$ws.Range("AB6").Value = 19
$ws.Range("AE6").Value = 7.2
$ws.Range("AF6").Value = 15
$ws.Range("AG6").Value = 12.5
$ws.Range("AH6").Value = 45
$ws.Range("AI6").Value = 40
